$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 12.58410635944461
$ws.Range("C2").Value = 7.246377296395679
$ws.Range("D2").Value = 7.804520551221181
$ws.Range("E2").Value = 13.01364591543346
$ws.Range("F2").Value = 38.8146951746339
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 29.36583346004858
$ws.Range("J2").Value = 10.3624905419177
$ws.Range("K2").Value = 10.06824227607603
$ws.Range("L2").Value = 11.13782472131582
$ws.Range("M2").Value = 14.98491546710325
$ws.Range("O2").Value = 30.07246463931372
$ws.Range("B3").Value = 12.38499528303822
$ws.Range("C3").Value = 7.194066479136485
$ws.Range("D3").Value = 7.793294854945007
$ws.Range("E3").Value = 13.03679928334488
$ws.Range("F3").Value = 38.92015916069266
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 29.46717665757866
$ws.Range("J3").Value = 10.37960624037694
$ws.Range("K3").Value = 9.922934195092827
$ws.Range("L3").Value = 11.14571950887442
$ws.Range("M3").Value = 14.95777706981444
$ws.Range("O3").Value = 30.16798427002439
$ws.Range("B4").Value = 12.26296533776107
$ws.Range("C4").Value = 7.161339742216431
$ws.Range("D4").Value = 7.787328045980961
$ws.Range("E4").Value = 13.05225367842496
$ws.Range("F4").Value = 38.9918714873321
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 29.53377498106243
$ws.Range("J4").Value = 10.39067825499223
$ws.Range("K4").Value = 9.833963979635687
$ws.Range("L4").Value = 11.15170351253038
$ws.Range("M4").Value = 14.94280531527304
$ws.Range("O4").Value = 30.23155128500857
$ws.Range("B5").Value = 12.21335399516977
$ws.Range("C5").Value = 7.147855362866677
$ws.Range("D5").Value = 7.785131243769025
$ws.Range("E5").Value = 13.05886336288557
$ws.Range("F5").Value = 39.02284296610626
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 29.56201461733676
$ws.Range("J5").Value = 10.39533214014762
$ws.Range("K5").Value = 9.797811637268685
$ws.Range("L5").Value = 11.15442854263819
$ws.Range("M5").Value = 14.93713375477335
$ws.Range("O5").Value = 30.25869163816903
$ws.Range("B6").Value = 12.20512489676291
$ws.Range("C6").Value = 7.14560749684972
$ws.Range("D6").Value = 7.784780696690831
$ws.Range("E6").Value = 13.0599797497145
$ws.Range("F6").Value = 39.02809129929332
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 29.56677025970038
$ws.Range("J6").Value = 10.39611349945196
$ws.Range("K6").Value = 9.791816034705828
$ws.Range("L6").Value = 11.15489835550091
$ws.Range("M6").Value = 14.93621805980369
$ws.Range("O6").Value = 30.26327293277623
$ws.Range("B7").Value = 12.26229570876982
$ws.Range("C7").Value = 7.161158479870198
$ws.Range("D7").Value = 7.787297466310051
$ws.Range("E7").Value = 13.05234155536017
$ws.Range("F7").Value = 38.99228210321875
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 29.53415137487424
$ws.Range("J7").Value = 10.39074044363255
$ws.Range("K7").Value = 9.833475942193711
$ws.Range("L7").Value = 11.15173910230878
$ws.Range("M7").Value = 14.94272708194971
$ws.Range("O7").Value = 30.23191230404968
$ws.Range("B8").Value = 12.515442126783
$ws.Range("C8").Value = 7.228469540969915
$ws.Range("D8").Value = 7.800458962107554
$ws.Range("E8").Value = 13.02137255627092
$ws.Range("F8").Value = 38.84961449975903
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 29.39986930230215
$ws.Range("J8").Value = 10.36827544588721
$ws.Range("K8").Value = 10.0181127915142
$ws.Range("L8").Value = 11.14031137705546
$ws.Range("M8").Value = 14.97520973029132
$ws.Range("O8").Value = 30.10437900006779
$ws.Range("B9").Value = 13.0110361157133
$ws.Range("C9").Value = 7.355435849691197
$ws.Range("D9").Value = 7.833526182768928
$ws.Range("E9").Value = 12.9704418864653
$ws.Range("F9").Value = 38.62509793578187
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 29.17121445272189
$ws.Range("J9").Value = 10.3286696768517
$ws.Range("K9").Value = 10.38037951472745
$ws.Range("L9").Value = 11.12688851664159
$ws.Range("M9").Value = 15.05213288979557
$ws.Range("O9").Value = 29.89331760710099
$ws.Range("B10").Value = 13.37122926788367
$ws.Range("C10").Value = 7.4453774665055
$ws.Range("D10").Value = 7.86212066768973
$ws.Range("E10").Value = 12.9389636937074
$ws.Range("F10").Value = 38.49389752951039
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 29.02432052228513
$ws.Range("J10").Value = 10.30225714405718
$ws.Range("K10").Value = 10.64430490157418
$ws.Range("L10").Value = 11.1224624267025
$ws.Range("M10").Value = 15.1164368620733
$ws.Range("O10").Value = 29.76206111843597
$ws.Range("B11").Value = 13.53354713994561
$ws.Range("C11").Value = 7.485511803575077
$ws.Range("D11").Value = 7.876033870329921
$ws.Range("E11").Value = 12.92592621830805
$ws.Range("F11").Value = 38.44155171748028
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 28.96206842909637
$ws.Range("J11").Value = 10.29081911561181
$ws.Range("K11").Value = 10.7634059408459
$ws.Range("L11").Value = 11.12161933640091
$ws.Range("M11").Value = 15.14731861439833
$ws.Range("O11").Value = 29.70752257859958
$ws.Range("B12").Value = 13.5947368301244
$ws.Range("C12").Value = 7.500592591093877
$ws.Range("D12").Value = 7.881429985092463
$ws.Range("E12").Value = 12.9211730505913
$ws.Range("F12").Value = 38.42278562840853
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 28.93915198385195
$ws.Range("J12").Value = 10.28657041556398
$ws.Range("K12").Value = 10.80832979856551
$ws.Range("L12").Value = 11.12146746386476
$ws.Range("M12").Value = 15.15924131839281
$ws.Range("O12").Value = 29.68761411169555
$ws.Range("B13").Value = 13.58157169738716
$ws.Range("C13").Value = 7.497349965159315
$ws.Range("D13").Value = 7.880262209591805
$ws.Range("E13").Value = 12.92218856232849
$ws.Range("F13").Value = 38.4267802563387
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 28.9440582297043
$ws.Range("J13").Value = 10.28748177974798
$ws.Range("K13").Value = 10.79866312633426
$ws.Range("L13").Value = 11.12149274265168
$ws.Range("M13").Value = 15.15666348598013
$ws.Range("O13").Value = 29.69186865386857
$ws.Range("B14").Value = 13.53858710596824
$ws.Range("C14").Value = 7.486754886398812
$ws.Range("D14").Value = 7.876475272844258
$ws.Range("E14").Value = 12.92553149091073
$ws.Range("F14").Value = 38.43998664871517
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 28.96016991166611
$ws.Range("J14").Value = 10.29046791847653
$ws.Range("K14").Value = 10.76710562213044
$ws.Range("L14").Value = 11.12160349241735
$ws.Range("M14").Value = 15.14829495826436
$ws.Range("O14").Value = 29.70586978119067
$ws.Range("B15").Value = 13.51222017285071
$ws.Range("C15").Value = 7.480249683673816
$ws.Range("D15").Value = 7.874172184889155
$ws.Range("E15").Value = 12.9276030570892
$ws.Range("F15").Value = 38.44821351043193
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 28.97012434695578
$ws.Range("J15").Value = 10.29230776523485
$ws.Range("K15").Value = 10.74775153795851
$ws.Range("L15").Value = 11.12169310094372
$ws.Range("M15").Value = 15.14319856695778
$ws.Range("O15").Value = 29.7145427950787
$ws.Range("B16").Value = 13.36058569187143
$ws.Range("C16").Value = 7.442738459637741
$ws.Range("D16").Value = 7.86122939738753
$ws.Range("E16").Value = 12.93984147800491
$ws.Range("F16").Value = 38.49746619423954
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 29.02848080036524
$ws.Range("J16").Value = 10.3030162290856
$ws.Range("K16").Value = 10.63649867866338
$ws.Range("L16").Value = 11.12254099853924
$ws.Range("M16").Value = 15.11445096060719
$ws.Range("O16").Value = 29.76572943584393
$ws.Range("B17").Value = 13.26712910506104
$ws.Range("C17").Value = 7.419523057965723
$ws.Range("D17").Value = 7.853519424105951
$ws.Range("E17").Value = 12.94767736061902
$ws.Range("F17").Value = 38.52956126717078
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 29.0654511170964
$ws.Range("J17").Value = 10.30973309162716
$ws.Range("K17").Value = 10.56797453476489
$ws.Range("L17").Value = 11.12336028020554
$ws.Range("M17").Value = 15.0972284934209
$ws.Range("O17").Value = 29.79845552319797
$ws.Range("B18").Value = 13.21323398082151
$ws.Range("C18").Value = 7.40609699616836
$ws.Range("D18").Value = 7.849170244135246
$ws.Range("E18").Value = 12.95230506323764
$ws.Range("F18").Value = 38.54871213701229
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 29.08714569049115
$ws.Range("J18").Value = 10.31365080029094
$ws.Range("K18").Value = 10.52847315778869
$ws.Range("L18").Value = 11.12394172621656
$ws.Range("M18").Value = 15.08747633252382
$ws.Range("O18").Value = 29.81776536156988
$ws.Range("B19").Value = 13.19496347420868
$ws.Range("C19").Value = 7.40153874689898
$ws.Range("D19").Value = 7.847712433613665
$ws.Range("E19").Value = 12.95389267183173
$ws.Range("F19").Value = 38.55531487752587
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 29.09456499797192
$ws.Range("J19").Value = 10.31498661437869
$ws.Range("K19").Value = 10.51508475272712
$ws.Range("L19").Value = 11.12415755070304
$ws.Range("M19").Value = 15.08420099993919
$ws.Range("O19").Value = 29.82438690785338
$ws.Range("B20").Value = 13.2770927604192
$ws.Range("C20").Value = 7.422001986891765
$ws.Range("D20").Value = 7.85433134539605
$ws.Range("E20").Value = 12.94683072827741
$ws.Range("F20").Value = 38.52607320651549
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 29.06147103851788
$ws.Range("J20").Value = 10.30901244807102
$ws.Range("K20").Value = 10.57527845270821
$ws.Range("L20").Value = 11.12326166569457
$ws.Range("M20").Value = 15.0990459841034
$ws.Range("O20").Value = 29.79492140389928
$ws.Range("B21").Value = 13.55122066046296
$ws.Range("C21").Value = 7.489870138270128
$ws.Range("D21").Value = 7.877584150175686
$ws.Range("E21").Value = 12.92454460656162
$ws.Range("F21").Value = 38.43607894126222
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 28.95541968878522
$ws.Range("J21").Value = 10.28958857738178
$ws.Range("K21").Value = 10.77637993736911
$ws.Range("L21").Value = 11.12156642713121
$ws.Range("M21").Value = 15.15074684538825
$ws.Range("O21").Value = 29.70173711188913
$ws.Range("B22").Value = 13.72874294267774
$ws.Range("C22").Value = 7.533540431564778
$ws.Range("D22").Value = 7.893522957543901
$ws.Range("E22").Value = 12.91105065817505
$ws.Range("F22").Value = 38.38341875953839
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 28.88993884420126
$ws.Range("J22").Value = 10.27737543802363
$ws.Range("K22").Value = 10.90676190826747
$ws.Range("L22").Value = 11.12143370158037
$ws.Range("M22").Value = 15.18586508164667
$ws.Range("O22").Value = 29.64517318261297
$ws.Range("B23").Value = 13.63416344596826
$ws.Range("C23").Value = 7.510297116177508
$ws.Range("D23").Value = 7.884949170079435
$ws.Range("E23").Value = 12.91815478120834
$ws.Range("F23").Value = 38.41096092143827
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 28.92453682900485
$ws.Range("J23").Value = 10.28384988354917
$ws.Range("K23").Value = 10.83728315430304
$ws.Range("L23").Value = 11.12141561899278
$ws.Range("M23").Value = 15.1670022339391
$ws.Range("O23").Value = 29.67496536745226
$ws.Range("B24").Value = 13.27258870364843
$ws.Range("C24").Value = 7.420881509506441
$ws.Range("D24").Value = 7.853964015849617
$ws.Range("E24").Value = 12.94721310837453
$ws.Range("F24").Value = 38.52764798152094
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 29.06326906192906
$ws.Range("J24").Value = 10.30933807626439
$ws.Range("K24").Value = 10.571976678377
$ws.Range("L24").Value = 11.1233059052726
$ws.Range("M24").Value = 15.09822383118105
$ws.Range("O24").Value = 29.79651763683028
$ws.Range("B25").Value = 12.87741170337106
$ws.Range("C25").Value = 7.321655085584218
$ws.Range("D25").Value = 7.823816014186602
$ws.Range("E25").Value = 12.98317432365403
$ws.Range("F25").Value = 38.67991272774475
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 29.22936381494227
$ws.Range("J25").Value = 10.33891057543508
$ws.Range("K25").Value = 10.28260175768464
$ws.Range("L25").Value = 11.12956192883053
$ws.Range("M25").Value = 15.02993462842796
$ws.Range("O25").Value = 29.94623505826393
